$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) with corrected values
# per "Correcion a Diebold Mariano y revision de Cap1"

$ws.Range("C2").Value = 0.5739478350893729
$ws.Range("D2").Value = 0.5718249855010749

$ws.Range("C3").Value = 0.9176656321648641
$ws.Range("D3").Value = 0.368745342755421

$ws.Range("C4").Value = 1.424647619281205
$ws.Range("D4").Value = 0.1682911359581423

$ws.Range("C5").Value = 0.4212650398169749
$ws.Range("D5").Value = 0.6776488631679007

$ws.Range("C6").Value = 0.5127073571464054
$ws.Range("D6").Value = 0.6132657044393484

$ws.Range("C7").Value = 0.7188609340482601
$ws.Range("D7").Value = 0.4797929885980197

$ws.Range("C8").Value = -0.05844310531803317
$ws.Range("D8").Value = 0.9539232877863097

$ws.Range("C9").Value = 0.2081560037133113
$ws.Range("D9").Value = 0.8370221699977993

$ws.Range("C10").Value = -0.49366620018346
$ws.Range("D10").Value = 0.6264344125763706

$ws.Range("C11").Value = -0.6274056657128783
$ws.Range("D11").Value = 0.5368514259868142
